$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1362.8334
$ws.Range("I135").Value = 1544.5
$ws.Range("J135").Value = 999.5
$ws.Range("K135").Value = 13900.5
$ws.Range("L135").Value = 8995.5
$ws.Range("M135").Value = -11365.5
$ws.Range("N135").Value = -14065.5
$ws.Range("H137").Value = 4761.5
$ws.Range("I137").Value = 4265.5
$ws.Range("K137").Value = 12796.5
$ws.Range("M137").Value = -10246.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3006.4546
$ws.Range("I61").Value = 2957.2
$ws.Range("K61").Value = 2957.2
$ws.Range("M61").Value = -2745.2
$ws.Range("H74").Value = 2686.5217
$ws.Range("I74").Value = 2155.0557
$ws.Range("K74").Value = 2155.0557
$ws.Range("M74").Value = -1281.0557
$ws.Range("H77").Value = 2686.5217
$ws.Range("I77").Value = 2155.0557
$ws.Range("K77").Value = 10775.2785
$ws.Range("M77").Value = -6407.2785
$ws.Range("H94").Value = 9997
$ws.Range("J94").Value = 9997
$ws.Range("L94").Value = 9997
$ws.Range("N94").Value = -11799
$ws.Range("H102").Value = 2146.5833
$ws.Range("I102").Value = 1809.8334
$ws.Range("J102").Value = 2483.3333
$ws.Range("K102").Value = 1809.8334
$ws.Range("L102").Value = 2483.3333
$ws.Range("M102").Value = -187.8334
$ws.Range("N102").Value = -5727.3333
$ws.Range("H105").Value = 12370
$ws.Range("J105").Value = 12370
$ws.Range("L105").Value = 12370
$ws.Range("N105").Value = -19358
$ws.Range("H122").Value = 3320.0588
$ws.Range("J122").Value = 3492.8572
$ws.Range("L122").Value = 10478.5716
$ws.Range("N122").Value = -15378.5716
$ws.Range("H125").Value = 87500
$ws.Range("J125").Value = 87500
$ws.Range("L125").Value = 87500
$ws.Range("N125").Value = -97340
$ws.Range("H136").Value = 3006.4546
$ws.Range("I136").Value = 2957.2
$ws.Range("K136").Value = 8871.599999999999
$ws.Range("M136").Value = -6321.599999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1688229.5
$ws.Range("I22").Value = 2530345
$ws.Range("J22").Value = 3998.5
$ws.Range("K22").Value = 2530345
$ws.Range("L22").Value = 3998.5
$ws.Range("M22").Value = -2530172
$ws.Range("N22").Value = -4344.5
$ws.Range("H86").Value = 1428.4286
$ws.Range("I86").Value = 1224.75
$ws.Range("J86").Value = 1700
$ws.Range("K86").Value = 1224.75
$ws.Range("L86").Value = 1700
$ws.Range("M86").Value = -101.75
$ws.Range("N86").Value = -3946
$ws.Range("H89").Value = 1428.4286
$ws.Range("I89").Value = 1224.75
$ws.Range("J89").Value = 1700
$ws.Range("K89").Value = 6123.75
$ws.Range("L89").Value = 8500
$ws.Range("M89").Value = -507.75
$ws.Range("N89").Value = -19732

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1255
$ws.Range("I16").Value = 1138.6666
$ws.Range("K16").Value = 1138.6666
$ws.Range("M16").Value = -851.6666
$ws.Range("H21").Value = 175
$ws.Range("J21").Value = 175
$ws.Range("L21").Value = 175
$ws.Range("N21").Value = -645
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1000
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1700
$ws.Range("H62").Value = 3259.6
$ws.Range("I62").Value = 2824.5
$ws.Range("K62").Value = 2824.5
$ws.Range("M62").Value = -2200.5
$ws.Range("H65").Value = 3259.6
$ws.Range("I65").Value = 2824.5
$ws.Range("K65").Value = 14122.5
$ws.Range("M65").Value = -11002.5
$ws.Range("H68").Value = 42500
$ws.Range("J68").Value = 50000
$ws.Range("L68").Value = 50000
$ws.Range("N68").Value = -51498
$ws.Range("H71").Value = 42500
$ws.Range("J71").Value = 50000
$ws.Range("L71").Value = 150000
$ws.Range("N71").Value = -157488
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H113").Value = 1255
$ws.Range("I113").Value = 1138.6666
$ws.Range("K113").Value = 1138.6666
$ws.Range("M113").Value = 1031.3334
$ws.Range("H132").Value = 1712.25
$ws.Range("I132").Value = 1283
$ws.Range("K132").Value = 3849
$ws.Range("M132").Value = -1319
$ws.Range("H134").Value = 2579.875
$ws.Range("I134").Value = 2650.1428
$ws.Range("J134").Value = 2088
$ws.Range("K134").Value = 7950.428400000001
$ws.Range("L134").Value = 6264
$ws.Range("M134").Value = -5415.428400000001
$ws.Range("N134").Value = -11334

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1000
$ws.Range("J50").Value = 1000
$ws.Range("L50").Value = 3000
$ws.Range("N50").Value = -3962
$ws.Range("H53").Value = 1000
$ws.Range("J53").Value = 1000
$ws.Range("L53").Value = 3000
$ws.Range("N53").Value = -3962
$ws.Range("H55").Value = 1171.25
$ws.Range("I55").Value = 761.6667
$ws.Range("K55").Value = 2285.0001
$ws.Range("M55").Value = -2108.0001
$ws.Range("H131").Value = 1292.1072
$ws.Range("I131").Value = 1695
$ws.Range("J131").Value = 1224.9584
$ws.Range("K131").Value = 5085
$ws.Range("L131").Value = 3674.8752
$ws.Range("M131").Value = -45
$ws.Range("N131").Value = -13754.8752
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 20015
$ws.Range("I35").Value = 20015
$ws.Range("K35").Value = 20015
$ws.Range("M35").Value = -19717
$ws.Range("H98").Value = 4060.5
$ws.Range("J98").Value = 4060.5
$ws.Range("L98").Value = 4060.5
$ws.Range("N98").Value = -10050.5
$ws.Range("H99").Value = 2999
$ws.Range("I99").Value = 2999
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 2999
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -753
$ws.Range("N99").ClearContents()
$ws.Range("H132").Value = 3468.3333
$ws.Range("I132").Value = 2912.1
$ws.Range("K132").Value = 8736.3
$ws.Range("M132").Value = -6206.299999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 2000
$ws.Range("I29").Value = 2000
$ws.Range("K29").Value = 2000
$ws.Range("M29").Value = -1705
$ws.Range("H64").Value = 21999.5
$ws.Range("J64").Value = 21999.5
$ws.Range("L64").Value = 21999.5
$ws.Range("N64").Value = -22449.5
$ws.Range("H67").Value = 21999.5
$ws.Range("J67").Value = 21999.5
$ws.Range("L67").Value = 21999.5
$ws.Range("N67").Value = -23559.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 8000
$ws.Range("J63").Value = 8000
$ws.Range("L63").Value = 8000
$ws.Range("N63").Value = -9248
$ws.Range("H66").Value = 8000
$ws.Range("J66").Value = 8000
$ws.Range("L66").Value = 24000
$ws.Range("N66").Value = -30240
$ws.Range("H136").Value = 1731.4286
$ws.Range("I136").Value = 1731.4286
$ws.Range("K136").Value = 5194.2858
$ws.Range("M136").Value = -2644.2858

